$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$times = @(
    "2021-10-05 13:40:36.041436",
    "2021-10-05 13:40:36.041449",
    "2021-10-05 13:40:36.041453",
    "2021-10-05 13:40:36.041456",
    "2021-10-05 13:40:36.041460",
    "2021-10-05 13:40:36.041463",
    "2021-10-05 13:40:36.041466",
    "2021-10-05 13:40:36.041469",
    "2021-10-05 13:40:36.041472",
    "2021-10-05 13:40:36.041475",
    "2021-10-05 13:40:36.041478",
    "2021-10-05 13:40:36.041481",
    "2021-10-05 13:40:36.041484",
    "2021-10-05 13:40:36.041487",
    "2021-10-05 13:40:36.041490",
    "2021-10-05 13:40:36.041493",
    "2021-10-05 13:40:36.041497",
    "2021-10-05 13:40:36.041500",
    "2021-10-05 13:40:36.041503",
    "2021-10-05 13:40:36.041506"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
